# Updates crypto price/volume data per the scraped cryptos.xlsx diff.
# Rows 48/49 additionally swap Coin name, Link, Price and Volume (ranking reorder).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.510.58"
$ws.Range("E2").Value = "  -0.67%  "
# Row 3
$ws.Range("D3").Value = "3.526.33"
$ws.Range("E3").Value = "  -0.79%  "
# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.02%  "
# Row 5
$ws.Range("D5").Value = "'608.49"
$ws.Range("E5").Value = "  -2.27%  "
# Row 6
$ws.Range("D6").Value = "'150.96"
$ws.Range("E6").Value = "  -2.74%  "
# Row 7
$ws.Range("D7").Value = "3.525.81"
$ws.Range("E7").Value = "  -0.59%  "
# Row 8
$ws.Range("E8").Value = "  -0.25%  "
# Row 9
$ws.Range("D9").Value = "'0.483"
$ws.Range("E9").Value = "  -0.95%  "
# Row 10
$ws.Range("E10").Value = "  -1.28%  "
# Row 11
$ws.Range("D11").Value = "'7.04"
$ws.Range("E11").Value = "  +1.19%  "
# Row 12
$ws.Range("D12").Value = "'0.425"
$ws.Range("E12").Value = "  -1.83%  "
# Row 13
$ws.Range("E13").Value = "  -2.60%  "
# Row 14
$ws.Range("D14").Value = "4.122.65"
$ws.Range("E14").Value = "  -0.91%  "
# Row 15
$ws.Range("D15").Value = "'31.99"
$ws.Range("E15").Value = "  -0.54%  "
# Row 16
$ws.Range("D16").Value = "3.527.44"
$ws.Range("E16").Value = "  -0.78%  "
# Row 17
$ws.Range("D17").Value = "67.515.85"
$ws.Range("E17").Value = "  -0.79%  "
# Row 18
$ws.Range("E18").Value = "  +0.03%  "
# Row 19
$ws.Range("D19").Value = "'6.43"
$ws.Range("E19").Value = "  +0.06%  "
# Row 20
$ws.Range("D20").Value = "'15.21"
$ws.Range("E20").Value = "  -2.76%  "
# Row 21
$ws.Range("D21").Value = "'446.35"
$ws.Range("E21").Value = "  -3.10%  "
# Row 22
$ws.Range("D22").Value = "'9.34"
$ws.Range("E22").Value = "  -4.24%  "
# Row 23
$ws.Range("D23").Value = "'0.624"
$ws.Range("E23").Value = "  -2.95%  "
# Row 24
$ws.Range("D24").Value = "'77.37"
$ws.Range("E24").Value = "  -0.67%  "
# Row 25
$ws.Range("E25").Value = "  +10.78%  "
# Row 26
$ws.Range("D26").Value = "3.666.14"
$ws.Range("E26").Value = "  -1.02%  "
# Row 27
$ws.Range("E27").Value = "  +0.22%  "
# Row 28
$ws.Range("D28").Value = "'10.20"
$ws.Range("E28").Value = "  -4.62%  "
# Row 29
$ws.Range("D29").Value = "'8.34"
$ws.Range("E29").Value = "  -0.54%  "
# Row 30
$ws.Range("D30").Value = "'2.50"
$ws.Range("E30").Value = "  -3.38%  "
# Row 31
$ws.Range("D31").Value = "'1.57"
$ws.Range("E31").Value = "  -3.88%  "
# Row 32
$ws.Range("E32").Value = "  +0.08%  "
# Row 33
$ws.Range("D33").Value = "'0.165"
$ws.Range("E33").Value = "  +4.52%  "
# Row 34
$ws.Range("D34").Value = "'25.76"
$ws.Range("E34").Value = "  -0.99%  "
# Row 35
$ws.Range("D35").Value = "'6.13"
$ws.Range("E35").Value = "  -0.99%  "
# Row 36
$ws.Range("D36").Value = "3.517.98"
$ws.Range("E36").Value = "  -1.11%  "
# Row 37
$ws.Range("E37").Value = "  -3.59%  "
# Row 38
$ws.Range("D38").Value = "'8.06"
$ws.Range("E38").Value = "  -0.45%  "
# Row 39
$ws.Range("E39").Value = "  -0.05%  "
# Row 40
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  -0.02%  "
# Row 41
$ws.Range("D41").Value = "'177.39"
$ws.Range("E41").Value = "  -0.41%  "
# Row 42
$ws.Range("D42").Value = "'2.17"
$ws.Range("E42").Value = "  +2.53%  "
# Row 43
$ws.Range("D43").Value = "'0.0876"
$ws.Range("E43").Value = "  -1.07%  "
# Row 44
$ws.Range("E44").Value = "  -3.58%  "
# Row 45
$ws.Range("D45").Value = "'0.881"
$ws.Range("E45").Value = "  -1.46%  "
# Row 46
$ws.Range("D46").Value = "'45.48"
$ws.Range("E46").Value = "  -0.93%  "
# Row 47
$ws.Range("E47").Value = "  +4.66%  "
# Row 48
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "'2.61"
$ws.Range("E48").Value = "  +1.07%  "
# Row 49
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'26.98"
$ws.Range("E49").Value = "  -5.32%  "
# Row 50
$ws.Range("D50").Value = "'7.59"
$ws.Range("E50").Value = "  -1.86%  "
# Row 51
$ws.Range("D51").Value = "'0.995"
$ws.Range("E51").Value = "  -1.60%  "
